# [#4200] Kickstart report, change table style and results log font size.
# Adds three new styles to the document's style sheet:
#   - "Plain Table 1" (styleId PlainTable1)  - built-in table style used for
#     the report tables (subtle grey borders, banded rows/cols, bold
#     header/first-row/first-col/last-row/last-col).
#   - "Table Common" (styleId TableCommon)   - custom table style carrying
#     the common cell margins used across the report's tables.
#   - "Normal Smaller" (styleId NormalSmaller) - custom paragraph style used
#     to shrink the results-log font size.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Plain Table 1 -- wdStyleTypeTable = 3
# ---------------------------------------------------------------------
$plainTable1 = $d.Styles.Add("Plain Table 1", 3)
$plainTable1.BaseStyle = "TableNormal"
$plainTable1.Priority = 99

# ---------------------------------------------------------------------
# 2) Table Common -- custom table style, wdStyleTypeTable = 3
# ---------------------------------------------------------------------
$tableCommon = $d.Styles.Add("Table Common", 3)
$tableCommon.BaseStyle = "TableNormal"
$tableCommon.Priority = 99

# ---------------------------------------------------------------------
# 3) Normal Smaller -- custom paragraph style, wdStyleTypeParagraph = 1
#    Used to shrink the results log font (size 10pt / w:sz 20 half-points,
#    black colour, Calibri for east-asian/complex-script runs).
# ---------------------------------------------------------------------
$normalSmaller = $d.Styles.Add("Normal Smaller", 1)
$normalSmaller.BaseStyle = "Normal"
$normalSmaller.QuickStyle = $true

$smallerFont = $normalSmaller.Font
$smallerFont.NameFarEast = "Calibri"
$smallerFont.NameBi = "Calibri"
$smallerFont.Size = 10
$smallerFont.Color = 0
